# Add a new "Swiss" test-data sheet (Switzerland market), modelled on the
# existing per-country sheets, and nudge a couple of pre-existing
# selections/active-tab state to match the saved workbook view.

$wb = $excel.ActiveWorkbook

# --- 1) Create the new "Swiss" sheet as a copy of the last existing
#        per-country sheet ("Czech"), placed after it (i.e. at the end). ---
$source = $wb.Worksheets.Item("Czech")
$source.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$swiss = $wb.Worksheets.Item($wb.Worksheets.Count)
$swiss.Name = "Swiss"

# --- 2) Fill in the Switzerland-specific market name / ticket code. ---
$swiss.Range("B2").Value = "Switzerland Market"
$swiss.Range("B4").Value = "NGC-3476/T2642/T2643/T2644"

# --- 3) Germany sheet: selection grew from A8 to A8:A16. ---
$germany = $wb.Worksheets.Item("Germany")
$germany.Range("A8:A16").Select()

# --- 4) Belgium sheet: it's no longer the active tab, and the whole sheet
#        (A1:XFD1048576) ends up selected. ---
$belgium = $wb.Worksheets.Item("Belgium")
$belgium.Cells.Select()

# --- 5) Swiss becomes the active tab, with B4 selected. ---
$swiss.Select()
$swiss.Range("B4").Select()
